$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 6.262800000000002
$ws.Range("B14").Value = 5.673200000000003
$ws.Range("C15").Value = -13.48319999999999
$ws.Range("B16").Value = 5.054499999999998
$ws.Range("B21").Value = 9.079499999999998
$ws.Range("C21").Value = -11.95670000000001
$ws.Range("C22").Value = -11.84920000000001
$ws.Range("B23").Value = 8.945500000000001
$ws.Range("C24").Value = -13.60149999999999
$ws.Range("B25").Value = 5.4453
$ws.Range("B26").Value = 5.061200000000007
$ws.Range("C27").Value = -12.25829999999999
$ws.Range("C28").Value = -13.9436
$ws.Range("B29").Value = 5.142900000000001
$ws.Range("C36").Value = -13.45780000000001
$ws.Range("C39").Value = -13.54560000000001
$ws.Range("B40").Value = 9.0677
$ws.Range("C45").Value = -13.38249999999999
$ws.Range("C48").Value = -11.91449999999999
$ws.Range("C49").Value = -13.90509999999999
$ws.Range("C52").Value = -10.81779999999999
$ws.Range("B53").Value = 5.802
$ws.Range("C53").Value = -10.9721
$ws.Range("C54").Value = -13.35299999999999
$ws.Range("B57").Value = 4.941399999999998
$ws.Range("C57").Value = -13.56079999999999
$ws.Range("B59").Value = 4.705699999999998
$ws.Range("B65").Value = 5.8975
$ws.Range("B69").Value = 5.593399999999995
$ws.Range("C70").Value = -12.37869999999999
$ws.Range("C71").Value = -11.1872
$ws.Range("B79").Value = 9.245500000000003
$ws.Range("B83").Value = 5.903999999999999
$ws.Range("C86").Value = -13.59509999999999
$ws.Range("C87").Value = -12.4921
$ws.Range("C89").Value = -13.08339999999999
$ws.Range("B91").Value = 6.246500000000003
$ws.Range("B93").Value = 5.521799999999997
$ws.Range("B100").Value = 4.931800000000004
$ws.Range("C101").Value = -13.4079
$ws.Range("B103").Value = 5.028000000000003
